$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; this shifts existing rows 50-70 down to 51-71,
# carrying over their formatting (including the date style on column D).
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new data record.
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 44603
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112001
$ws.Cells.Item(50, 7).Value = "Berenjena"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 170
$ws.Cells.Item(50, 11).Value = 9000
$ws.Cells.Item(50, 12).Value = 9500
$ws.Cells.Item(50, 13).Value = 9235
$ws.Cells.Item(50, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 154
$ws.Cells.Item(50, 17).Value = 60
$ws.Cells.Item(50, 18).Value = "Hortaliza"
